$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target data for the country table (rows 4-188), already sorted by
# "Casos totales" (column B) descending as of the update timestamp below.
# Columns: Name, Casos totales, Nuevos casos, Casos activos, Recuperados,
#          Casos criticos, Muertes hoy, Muertes
$data = @"
China	81008	41	71740	6013	1927	7	3255
Italia	47021	0	5129	37860	2655	0	4032
España	21571	0	1588	18890	939	0	1093
Alemania	19848	0	180	19600	2	0	68
Estados Unidos	19773	390	147	19351	64	19	275
Iran	19644	0	6745	11466	0	0	1433
Francia	12612	0	1587	10575	1297	0	450
Corea del Sur	8799	147	2612	6085	59	8	102
Suiza	5615	0	15	5544	141	0	56
Reino Unido	3983	0	65	3741	20	0	177
Paises Bajos	2994	0	2	2886	210	0	106
Austria	2649	0	9	2634	14	0	6
Belgica	2257	0	204	2016	164	0	37
Noruega	1959	0	1	1951	27	0	7
Suecia	1639	0	16	1607	21	0	16
Dinamarca	1255	0	1	1245	37	0	9
Canada	1087	0	14	1061	1	0	12
Australia	1051	123	46	998	2	0	7
Malasia	1030	0	87	940	26	0	3
Portugal	1020	0	5	1009	26	0	6
Japon	1007	0	215	757	50	0	35
Brasil	970	0	2	957	18	0	11
Chequia	833	0	4	829	6	0	0
Crucero	712	0	567	137	15	0	8
Israel	705	0	15	689	10	0	1
Irlanda	683	0	5	675	6	0	3
Turquia	670	0	0	661	0	0	9
Pakistan	501	0	13	485	0	0	3
Grecia	495	0	19	466	20	0	10
Luxemburgo	484	0	6	473	1	0	5
Catar	470	0	10	460	6	0	0
Finlandia	450	0	10	440	2	0	0
Chile	434	0	6	428	7	0	0
Ecuador	426	0	3	416	2	0	7
Polonia	425	0	13	407	3	0	5
Tailandia	411	89	42	368	1	0	1
Islandia	409	0	5	404	1	0	0
Singapur	385	0	131	252	14	2	2
Indonesia	369	0	17	320	0	0	32
Arabia Saudita	344	0	8	336	0	0	0
Eslovenia	341	0	0	340	9	0	1
Rumania	308	0	31	277	11	0	0
Barein	298	0	125	172	4	0	1
Egipto	285	0	42	235	0	0	8
Estonia	283	0	1	282	1	0	0
India	270	21	23	242	0	0	5
Peru	263	0	1	258	5	0	4
Filipinas	262	32	8	236	1	0	18
Hong Kong	256	0	98	154	4	0	4
Rusia	253	0	12	240	0	0	1
Irak	208	0	49	142	0	0	17
Mexico	203	39	4	197	1	1	2
Sudafrica	202	0	0	202	0	0	0
Panama	200	0	1	198	7	0	1
Libano	177	0	4	169	3	0	4
Kuwait	159	0	22	137	5	0	0
Colombia	158	13	1	157	0	0	0
Argentina	158	0	3	152	0	0	3
San Marino	151	0	4	133	12	0	14
Emiratos Arabes Unidos	140	0	31	107	2	0	2
Eslovaquia	137	0	0	137	2	0	0
Armenia	136	0	1	135	2	0	0
Serbia	135	0	2	132	4	0	1
Taiwan	135	0	28	105	0	0	2
Croacia	130	0	5	124	0	0	1
Bulgaria	127	0	3	121	0	0	3
Costa Rica	113	0	2	109	2	0	2
Letonia	111	0	1	110	0	0	0
Uruguay	110	0	0	110	0	0	0
Argelia	94	0	32	51	0	0	11
Vietnam	91	0	17	74	0	0	0
Bosnia y Herzegovina	89	0	2	87	1	0	0
Marruecos	86	0	2	81	1	0	3
Jordania	85	0	1	84	0	0	0
Hungria	85	0	7	74	6	0	4
Islas Feroe	80	0	3	77	0	0	0
Brunei	78	0	1	77	2	0	0
Republica de Macedonia	76	0	1	75	1	0	0
Republica de Chipre	75	0	0	75	1	0	0
Principado de Andorra	75	0	1	74	2	0	0
Sri Lanka	73	0	3	70	0	0	0
Republica Dominicana	72	0	0	70	0	0	2
Albania	70	0	2	66	2	0	2
Lituania	69	0	1	67	1	0	1
Bielorrusia	69	0	15	54	0	0	0
Moldavia	66	0	1	64	3	0	1
Venezuela	65	0	1	64	0	0	0
Malta	64	0	2	62	1	0	0
Tunez	54	0	1	52	7	0	1
Kazajistan	52	0	0	52	0	0	0
Nueva Zelanda	52	0	0	52	0	0	0
Camboya	51	0	1	50	0	0	0
Guadalupe	51	0	0	50	4	0	1
Oman	48	0	13	35	0	0	0
Estado de Palestina	48	0	17	31	0	0	0
Senegal	47	0	5	42	0	0	0
Georgia	44	0	1	43	1	0	0
Azerbaiyan	44	0	7	36	0	0	1
Ucrania	41	0	1	37	0	0	3
Burkina Faso	40	0	4	35	0	0	1
Reunion	38	0	0	38	0	0	0
Uzbekistan	33	0	0	33	0	0	0
Martinica	32	0	0	31	7	0	1
Liechtenstein	28	0	0	28	0	0	0
Camerun	27	0	2	25	0	0	0
Honduras	24	0	0	24	0	0	0
Afganistan	24	0	1	23	0	0	0
Cuba	21	0	0	20	0	0	1
Banglades	20	0	3	16	1	0	1
Jamaica	19	0	2	16	0	0	1
Consejo Danes para los Refugiados	18	0	0	18	0	0	0
Paraguay	18	0	0	18	1	0	0
Ruanda	17	0	0	17	0	0	0
Macao	17	0	10	7	0	0	0
Bolivia	16	0	0	16	0	0	0
Ghana	16	0	0	16	0	0	0
Guayana Francesa	15	0	0	15	0	0	0
Guyana	15	0	0	14	0	0	1
Guam	14	0	0	14	0	0	0
Montenegro	14	0	0	14	0	0	0
Puerto Rico	14	0	0	14	0	0	0
Maldivas	13	0	2	11	0	0	0
Mauricio	12	0	0	12	0	0	0
Guatemala	12	0	0	11	0	0	1
Nigeria	12	0	1	11	0	0	0
Monaco	11	0	0	11	0	0	0
Polinesia Francesa	11	0	0	11	0	0	0
Mongolia	10	4	0	10	0	0	0
Gibraltar	10	0	2	8	0	0	0
Trinidad yTobago	9	0	0	9	0	0	0
Togo	9	0	0	9	0	0	0
Etiopia	9	0	0	9	0	0	0
Costa de Marfil	9	0	1	8	0	0	0
Kenia	7	0	0	7	0	0	0
Seychelles	7	0	0	7	0	0	0
Mayotte	7	0	0	7	0	0	0
Barbados	6	0	0	6	0	0	0
Islas Virgenes de los Estados Unidos	6	3	0	6	0	0	0
Guinea Ecuatorial	6	0	0	6	0	0	0
Kirguistan	6	0	0	6	0	0	0
Tanzania	6	0	0	6	0	0	0
Aruba	5	0	1	4	0	0	0
Bahamas	4	0	0	4	0	0	0
Surinam	4	0	0	4	0	0	0
San Martin (Parte Francesa)	4	0	0	4	0	0	0
Gabon	4	0	0	3	0	0	1
El Salvador	3	2	0	3	0	0	0
Republica de Africa Central	3	0	0	3	0	0	0
Congo	3	0	0	3	0	0	0
Namibia	3	0	0	3	0	0	0
San Bartolome	3	0	0	3	0	0	0
Madagascar	3	0	0	3	0	0	0
Islas Caimanes	3	0	0	2	0	0	1
Curazao	3	0	0	2	0	0	1
Fiyi	2	1	0	2	0	0	0
Nueva Caledonia	2	0	0	2	0	0	0
Bermudas	2	0	0	2	0	0	0
Guinea	2	0	0	2	0	0	0
Liberia	2	0	0	2	0	0	0
Benin	2	0	0	2	0	0	0
Nicaragua	2	0	0	2	0	0	0
Santa Lucia	2	0	0	2	0	0	0
Mauritania	2	0	0	2	0	0	0
Groenlandia	2	0	0	2	0	0	0
Zambia	2	0	0	2	0	0	0
Isla de Man	2	0	0	2	0	0	0
Butan	2	0	0	2	0	0	0
Haiti	2	0	0	2	0	0	0
Sudan	2	0	0	1	0	0	1
Papua Nueva Guinea	1	0	0	1	0	0	0
Angola	1	0	0	1	0	0	0
Zimbabue	1	0	0	1	0	0	0
Santa Sede	1	0	0	1	0	0	0
Antigua y Barbuda	1	0	0	1	0	0	0
Gambia	1	0	0	1	0	0	0
Somalia	1	0	0	1	0	0	0
San Martin (Parte Holandesa)	1	0	0	1	0	0	0
Republica de Yibuti	1	0	0	1	0	0	0
Montserrat	1	0	0	1	0	0	0
San Vicente y las Granadinas	1	0	0	1	0	0	0
Cabo Verde	1	0	0	1	0	0	0
Suazilandia	1	0	0	1	0	0	0
Republica del Chad	1	0	0	1	0	0	0
Niger	1	0	0	1	0	0	0
Nepal	1	0	1	0	0	0	0
"@

$lines = $data -split "`n"
$row = 4
foreach ($line in $lines) {
    if ($line.Trim().Length -eq 0) { continue }
    $parts = $line -split "`t"
    $ws.Cells.Item($row, 1).Value = $parts[0]
    $ws.Cells.Item($row, 2).Value = [double]$parts[1]
    $ws.Cells.Item($row, 3).Value = [double]$parts[2]
    $ws.Cells.Item($row, 4).Value = [double]$parts[3]
    $ws.Cells.Item($row, 5).Value = [double]$parts[4]
    $ws.Cells.Item($row, 6).Value = [double]$parts[5]
    $ws.Cells.Item($row, 7).Value = [double]$parts[6]
    $ws.Cells.Item($row, 8).Value = [double]$parts[7]
    $row++
}

$ws.Range("A1").Value = "Datos actualizados a 21 de Marzo de 2020 a las 06:46"
